# Reclassify tile grid "map" sheet (A1:N14) using pixel-histogram labels.
# Letters: L=Land, B=Building/Blue, G=Green, M=Marker, W=Water
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("map")

# Row 1
$ws.Range("A1:N1").Value = "L"
# Row 2
$ws.Range("A2:N2").Value = "L"
# Row 3
$ws.Range("A3:L3").Value = "L"
$ws.Range("M3:N3").Value = "G"
# Row 4
$ws.Range("A4").Value = "B"
$ws.Range("B4:C4").Value = "G"
$ws.Range("D4:E4").Value = "L"
$ws.Range("F4").Value = "B"
$ws.Range("G4:M4").Value = "L"
$ws.Range("N4").Value = "W"
# Row 5
$ws.Range("A5:C5").Value = "L"
$ws.Range("D5").Value = "G"
$ws.Range("E5:I5").Value = "L"
$ws.Range("J5").Value = "G"
$ws.Range("K5:N5").Value = "L"
# Row 6
$ws.Range("A6:D6").Value = "L"
$ws.Range("E6:F6").Value = "G"
$ws.Range("G6:H6").Value = "L"
$ws.Range("I6").Value = "G"
$ws.Range("J6:K6").Value = "L"
$ws.Range("L6:M6").Value = "G"
$ws.Range("N6").Value = "L"
# Row 7
$ws.Range("A7:B7").Value = "L"
$ws.Range("C7").Value = "G"
$ws.Range("D7:G7").Value = "L"
$ws.Range("H7").Value = "B"
$ws.Range("I7").Value = "L"
$ws.Range("J7").Value = "G"
$ws.Range("K7:L7").Value = "L"
$ws.Range("M7:N7").Value = "G"
# Row 8
$ws.Range("A8:C8").Value = "G"
$ws.Range("D8").Value = "L"
$ws.Range("E8:F8").Value = "G"
$ws.Range("G8").Value = "M"
$ws.Range("H8").Value = "G"
$ws.Range("I8").Value = "L"
$ws.Range("J8:L8").Value = "G"
$ws.Range("M8:N8").Value = "L"
# Row 9
$ws.Range("A9").Value = "L"
$ws.Range("B9").Value = "M"
$ws.Range("C9").Value = "B"
$ws.Range("D9").Value = "L"
$ws.Range("E9").Value = "G"
$ws.Range("F9:N9").Value = "L"
# Row 10
$ws.Range("A10:J10").Value = "L"
$ws.Range("K10").Value = "G"
$ws.Range("L10:N10").Value = "L"
# Row 11
$ws.Range("A11:B11").Value = "L"
$ws.Range("C11").Value = "G"
$ws.Range("D11:N11").Value = "L"
# Row 12
$ws.Range("A12:I12").Value = "L"
$ws.Range("J12").Value = "W"
$ws.Range("K12").Value = "L"
$ws.Range("L12:M12").Value = "W"
$ws.Range("N12").Value = "B"
# Row 13
$ws.Range("A13:C13").Value = "L"
$ws.Range("D13").Value = "W"
$ws.Range("E13:G13").Value = "L"
$ws.Range("H13").Value = "W"
$ws.Range("I13:N13").Value = "L"
# Row 14
$ws.Range("A14:N14").Value = "L"

Write-Host "Reclassified $($ws.Range('A1:N14').Cells.Count) tiles on sheet $($ws.Name)"
